$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiment tracking")

# Fill the Experiment ID column down first (continuing the 001-004 pattern)
$ws.Range("A6").Value = "005"
$ws.Range("A7").Value = "006"
$ws.Range("A5").Copy()
$ws.Range("A6:A7").PasteSpecial(-4122)

# Row 6 - Experiment 005
$ws.Range("B6").Value = "005_fgfr2_ligand_x50.yaml"
$ws.Range("C6").Value = "FGFR2_v1"
$ws.Range("D6").Value = "ligand_v1"
$ws.Range("E6").Value = "A:378, A:398"
$ws.Range("F6").Value = "6.0 Å"
$ws.Range("G6").Value = "empty"
$ws.Range("H6").Value = "added 50 ligands"
$ws.Range("I6").Value = "005_fgfr2_ligand_x50_model.cif"
$ws.Range("J6").Value = "N/A"
$ws.Range("K6").Value = "Since the number of ligands is too large, the CIF file cannot be opened properly. "

# Row 7 - Experiment 006
$ws.Range("B7").Value = "006_fgfr2_ECD_ICD_torsion.yaml"
$ws.Range("C7").Value = "FGFR2_v1"
$ws.Range("D7").Value = "ligand_v1"
$ws.Range("E7").Value = "A:378, A:398"
$ws.Range("F7").Value = "6.0 Å"
$ws.Range("G7").Value = "empty"
$ws.Range("H7").Value = "Separate ECD and ICD via CXC file"
$ws.Range("I7").Value = "006_fgfr2_ECD_ICD_torsion.cif"
$ws.Range("J7").Value = "N/A"
$ws.Range("K7").Value = "The two domains appear to be successfully separated, although the ECD still shows interactions with the TM region."

$ws.Range("K10").Select()
